$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Salary" column (E) with header and values
$ws.Range("E1").Value = "Salary"
$ws.Range("E2").Value = 200000.0
$ws.Range("E3").Value = 190000.0
$ws.Range("E4").Value = 175000.0
$ws.Range("E5").Value = 210000.0

# Update selection to match the diff (F1:F5 selected, F1 active)
$ws.Range("F1:F5").Select()
